$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  ,@("D2", "67.308.71", $false)
  ,@("E2", "  +1.65%  ", $false)
  ,@("D3", "3.297.90", $false)
  ,@("E3", "  -0.89%  ", $false)
  ,@("D4", "1.00", $true)
  ,@("E4", "  +0.00%  ", $false)
  ,@("D5", "577.77", $true)
  ,@("E5", "  -0.92%  ", $false)
  ,@("D6", "176.35", $true)
  ,@("E6", "  -5.16%  ", $false)
  ,@("D7", "0.998", $true)
  ,@("E7", "  -0.17%  ", $false)
  ,@("D8", "0.586", $true)
  ,@("E8", "  +1.37%  ", $false)
  ,@("D9", "3.297.84", $false)
  ,@("E9", "  -0.78%  ", $false)
  ,@("D10", "0.176", $true)
  ,@("E10", "  -3.01%  ", $false)
  ,@("D11", "0.578", $true)
  ,@("E11", "  -0.16%  ", $false)
  ,@("D12", "45.74", $true)
  ,@("E12", "  -2.85%  ", $false)
  ,@("D13", "0.0000270", $true)
  ,@("E13", "  +0.62%  ", $false)
  ,@("D14", "688.34", $true)
  ,@("D15", "3.832.12", $false)
  ,@("E15", "  -0.73%  ", $false)
  ,@("D16", "8.38", $true)
  ,@("E16", "  -1.22%  ", $false)
  ,@("D17", "67.498.02", $false)
  ,@("E17", "  +1.80%  ", $false)
  ,@("E18", "  +1.35%  ", $false)
  ,@("D19", "3.297.93", $false)
  ,@("E19", "  -0.87%  ", $false)
  ,@("D20", "17.39", $true)
  ,@("E20", "  -3.08%  ", $false)
  ,@("D21", "10.80", $true)
  ,@("E21", "  -3.20%  ", $false)
  ,@("D22", "0.892", $true)
  ,@("E22", "  -0.45%  ", $false)
  ,@("D23", "17.13", $true)
  ,@("E23", "  -3.72%  ", $false)
  ,@("D24", "5.20", $true)
  ,@("E24", "  +3.10%  ", $false)
  ,@("D25", "99.56", $true)
  ,@("E25", "  -2.81%  ", $false)
  ,@("D26", "3.89", $true)
  ,@("E26", "  -2.55%  ", $false)
  ,@("D27", "2.71", $true)
  ,@("E27", "  -2.65%  ", $false)
  ,@("B28", "EthereumClassic", $false)
  ,@("C28", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", $false)
  ,@("D28", "33.95", $true)
  ,@("E28", "  +6.80%  ", $false)
  ,@("B29", "RenderToken", $false)
  ,@("C29", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", $false)
  ,@("D29", "9.29", $true)
  ,@("E29", "  -2.59%  ", $false)
  ,@("D30", "8.44", $true)
  ,@("E30", "  -0.57%  ", $false)
  ,@("D31", "6.72", $true)
  ,@("E31", "  -1.05%  ", $false)
  ,@("D32", "573.90", $true)
  ,@("E32", "  -4.40%  ", $false)
  ,@("D33", "3.887.38", $false)
  ,@("E33", "  +0.99%  ", $false)
  ,@("D34", "10.89", $true)
  ,@("E34", "  -1.06%  ", $false)
  ,@("D35", "0.103", $true)
  ,@("E35", "  -2.38%  ", $false)
  ,@("E36", "  +0.05%  ", $false)
  ,@("D37", "3.35", $true)
  ,@("E37", "  -14.37%  ", $false)
  ,@("D38", "55.42", $true)
  ,@("E38", "  -1.10%  ", $false)
  ,@("D39", "0.129", $true)
  ,@("E39", "  +1.65%  ", $false)
  ,@("D40", "3.40", $true)
  ,@("E40", "  -0.82%  ", $false)
  ,@("D41", "2.60", $true)
  ,@("E41", "  -3.76%  ", $false)
  ,@("D42", "31.91", $true)
  ,@("E42", "  -2.79%  ", $false)
  ,@("D43", "0.0₃0675", $false)
  ,@("E43", "  -3.65%  ", $false)
  ,@("D44", "0.329", $true)
  ,@("E44", "  -2.36%  ", $false)
  ,@("D45", "2.99", $true)
  ,@("E45", "  -5.40%  ", $false)
  ,@("D46", "0.0406", $true)
  ,@("E46", "  -1.74%  ", $false)
  ,@("D47", "0.128", $true)
  ,@("E47", "  +0.19%  ", $false)
  ,@("B48", "ThetaToken", $false)
  ,@("C48", "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta", $false)
  ,@("D48", "2.57", $true)
  ,@("E48", "  +0.60%  ", $false)
  ,@("B49", "FirstDigitalUSD", $false)
  ,@("C49", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", $false)
  ,@("D49", "1.01", $true)
  ,@("E49", "  +0.29%  ", $false)
  ,@("D50", "1.38", $true)
  ,@("E50", "  +5.38%  ", $false)
  ,@("D51", "130.34", $true)
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newVal = $u[1]
    $forceText = $u[2]
    $rng = $ws.Range($cellRef)
    if ($forceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $newVal
}
